# Commit: "Added CSV's and timemaps to src_data folder and timemaps folder"
#
# For this workbook the meaningful content change is that four
# Department rows (7, 9, 10 and 11) which previously showed the
# literal placeholder text "No URL" in column B (because no catalog
# hyperlink existed for that department) have had that placeholder
# text cleared out, leaving an empty-but-still-styled cell. All other
# "No URL" cells (rows 41-45) are left untouched, and none of these
# four rows ever had a real hyperlink attached, so no hyperlink
# relationships need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()

# Restore the cursor position to match the author's saved selection.
[void]$ws.Range("G12").Select()
